$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 (I0) and J1 (IF), formatted like the existing
# header cells (bold, centered, bordered) by copying H1's format.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-42 for new columns I and J
$iValues = @(7,8,10,6,7,7,6,6,7,7,6,8,7,7,8,7,6,7,7,7,6,6,7,8,8,7,7,8,9,8,7,8,5,7,7,5,7,6,7,5,6)
$jValues = @(7,8,10,7,7,7,7,6,7,8,7,8,7,8,8,7,6,8,7,7,6,6,7,8,8,8,7,8,9,8,7,8,6,7,7,5,7,6,7,5,6)

for ($r = 2; $r -le 42; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
